$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 7-14 in column L ("Cost ($)") become "Not Available" text entries
# instead of the numeric cost values read from the PM Meter import.
foreach ($row in 7..14) {
    $cell = $ws.Range("L$row")
    $cell.Value = "Not Available"
    $cell.NumberFormat = "@"
    $cell.HorizontalAlignment = 1
}
